$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8:10 (XBB.1.5 vs LY-CoV1404/AZD1061/AZD8895) are the template for the
# new "XBB.1.5_DJM" comparison rows 11:13 - same antibody/PDB/active-residue
# data, just a different Spike RBD variant + PDB file, and no HADDOCK results
# yet (columns K:AA stay empty). Copy formats row-by-row so each new row
# picks up the exact per-column style already used on its template row.
$ws.Range("A8:J8").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)
$ws.Range("A9:J9").Copy()
$ws.Range("A12:J12").PasteSpecial(-4122)
$ws.Range("A10:J10").Copy()
$ws.Range("A13:J13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value() = $ws.Range("A8").Value()
$ws.Range("B11").Value() = $ws.Range("B8").Value()
$ws.Range("C11").Value() = $ws.Range("C8").Value()
$ws.Range("D11").Value() = "XBB.1.5_DJM"
$ws.Range("E11").Value() = "XBB_1_5_DJM_ff8e4_relaxed_rank_1_model_1_renumbered.pdb"
$ws.Range("F11").Value() = $ws.Range("F8").Value()
$ws.Range("G11").Formula = "=_xlfn.CONCAT(SUBSTITUTE(D11,""."",""_""),""__"",A11)"
$ws.Range("H11").Value() = $ws.Range("H8").Value()
$ws.Range("J11").Value() = $true

$ws.Range("A12").Value() = $ws.Range("A9").Value()
$ws.Range("B12").Value() = $ws.Range("B9").Value()
$ws.Range("C12").Value() = $ws.Range("C9").Value()
$ws.Range("D12").Value() = "XBB.1.5_DJM"
$ws.Range("E12").Value() = "XBB_1_5_DJM_ff8e4_relaxed_rank_1_model_1_renumbered.pdb"
$ws.Range("F12").Value() = $ws.Range("F9").Value()
$ws.Range("G12").Formula = "=_xlfn.CONCAT(SUBSTITUTE(D12,""."",""_""),""__"",A12)"
$ws.Range("H12").Value() = $ws.Range("H9").Value()
$ws.Range("J12").Value() = $true

$ws.Range("A13").Value() = $ws.Range("A10").Value()
$ws.Range("B13").Value() = $ws.Range("B10").Value()
$ws.Range("C13").Value() = $ws.Range("C10").Value()
$ws.Range("D13").Value() = "XBB.1.5_DJM"
$ws.Range("E13").Value() = "XBB_1_5_DJM_ff8e4_relaxed_rank_1_model_1_renumbered.pdb"
$ws.Range("F13").Value() = $ws.Range("F10").Value()
$ws.Range("G13").Formula = "=_xlfn.CONCAT(SUBSTITUTE(D13,""."",""_""),""__"",A13)"
$ws.Range("H13").Value() = $ws.Range("H10").Value()
$ws.Range("J13").Value() = $true

# Match the row heights used by the other rows in the same antibody cycle
# (LY-CoV1404 / AZD1061 / AZD8895 -> 45 / 60 / 45).
$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(8).RowHeight
$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(9).RowHeight
$ws.Rows.Item(13).RowHeight = $ws.Rows.Item(10).RowHeight

# Extend the merged "RBD Active Residues" column to cover the new rows.
$ws.Range("I2:I10").UnMerge()
$ws.Range("I2:I13").Merge()
